$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "1143367910"
$ws.Range("D17").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E17").Value = "2008"
$ws.Range("F17").Value = 35112

$ws.Range("C18").Value = "1143367910"
$ws.Range("D18").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E18").Value = "2101"

$ws.Range("C19").Value = "1201219362"
$ws.Range("D19").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E19").Value = "2101"

$ws.Range("C20").Value = "1143367910"
$ws.Range("D20").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E20").Value = "2102"

$ws.Range("C21").Value = "1201219362"
$ws.Range("D21").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E21").Value = "2102"

$ws.Range("C22").Value = "1143367910"
$ws.Range("D22").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E22").Value = "2103"

$ws.Range("C23").Value = "1201219362"
$ws.Range("D23").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E23").Value = "2103"

$ws.Range("C24").Value = "1047424362"
$ws.Range("D24").Value = "VANESSA OSORIO SIMANCAS"
$ws.Range("E24").Value = "2103"

$ws.Range("C25").Value = "1143367910"
$ws.Range("D25").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E25").Value = "2104"
$ws.Range("F25").Value = 35112

$ws.Range("C26").Value = "1201219362"
$ws.Range("D26").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E26").Value = "2104"

$ws.Range("C27").Value = "1047424362"
$ws.Range("D27").Value = "VANESSA OSORIO SIMANCAS"
$ws.Range("E27").Value = "2104"

$ws.Range("C28").Value = "1143367910"
$ws.Range("D28").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E28").Value = "2105"

$ws.Range("C29").Value = "1201219362"
$ws.Range("D29").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E29").Value = "2105"

$ws.Range("C30").Value = "1047424362"
$ws.Range("D30").Value = "VANESSA OSORIO SIMANCAS"
$ws.Range("E30").Value = "2105"

$ws.Range("C31").Value = "1143367910"
$ws.Range("D31").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E31").Value = "2106"

$ws.Range("C32").Value = "1201219362"
$ws.Range("D32").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E32").Value = "2106"
$ws.Range("F32").Value = 35112

$ws.Range("C33").Value = "1047424362"
$ws.Range("D33").Value = "VANESSA OSORIO SIMANCAS"
$ws.Range("E33").Value = "2106"

$ws.Range("C34").Value = "1143367910"
$ws.Range("D34").Value = "JESUS DAVID LOPEZ BROME"
$ws.Range("E34").Value = "2107"
$ws.Range("F34").Value = 29260

$ws.Range("C35").Value = "1201219362"
$ws.Range("D35").Value = "OLISMAIDA GARCIA SOTO"
$ws.Range("E35").Value = "2107"
$ws.Range("F35").Value = 29260

$ws.Range("C36").Value = "1047424362"
$ws.Range("D36").Value = "VANESSA OSORIO SIMANCAS"
$ws.Range("E36").Value = "2107"
$ws.Range("F36").Value = 29260
